# Apply Wikipedia cabinet data refresh edits to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 - Chris Bishop
$ws.Range("H4").Value = "Bachelor of Arts in History and Politics, first-class honours in Law"
$ws.Range("K4").Value = "Arts and Law"

# Row 7 - Erica Stanford
$ws.Range("J7").Value = ""

# Row 9 - Louise Upston
$ws.Range("H9").Value = "Master of Business Administration"
$ws.Range("K9").Value = ""

# Row 10 - Judith Collins
$ws.Range("J10").Value = "1977, 1978, 1979, 2020"

# Row 14 - Matt Doocey
$ws.Range("J14").Value = "no specific years mentioned"
$ws.Range("K14").Value = ""

# Row 15 - Simon Watts
# Excel strips a single leading apostrophe used as a text-qualifier prefix,
# so prepend an extra apostrophe to end up with the literal text "''".
$ws.Range("J15").Value = "'''"

# Row 16 - David Seymour
# Force text formatting so the numeric-looking string isn't converted to a number.
$ws.Range("J16").NumberFormat = "@"
$ws.Range("J16").Value = "2014"

# Row 19 - Winston Peters
$ws.Range("I19").Value = "University of Auckland"
$ws.Range("K19").Value = "Arts, Science"

# Row 20 - Shane Jones
$ws.Range("I20").Value = "Victoria University of Wellington, Harvard Kennedy School at Harvard University"

# Row 26 - Andrew Bayly
$ws.Range("H26").Value = "Bachelor of Commerce in Accounting and Finance"
$ws.Range("J26").Value = ""

# Row 27 - Andrew Hoggard
$ws.Range("K27").Value = "Arts"

# Row 30 - Simon Court
$ws.Range("J30").Value = "not specified"
